# Update "想去人数" (F column) values across the worksheets to reflect
# freshly generated output (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 321
$ws1.Range("F7").Value = 867
$ws1.Range("F8").Value = 53
$ws1.Range("F9").Value = 517
$ws1.Range("F12").Value = 1144
$ws1.Range("F14").Value = 241
$ws1.Range("F15").Value = 35
$ws1.Range("F16").Value = 414
$ws1.Range("F17").Value = 6646
$ws1.Range("F23").Value = 37
$ws1.Range("F24").Value = 3399
$ws1.Range("F26").Value = 2107
$ws1.Range("F28").Value = 4517
$ws1.Range("F29").Value = 143
$ws1.Range("F34").Value = 1694
$ws1.Range("F36").Value = 170
$ws1.Range("F39").Value = 1208
$ws1.Range("F40").Value = 1783
$ws1.Range("F41").Value = 2138

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 81

# Sheet "本地生活"
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 1229

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1229
$ws4.Range("F7").Value = 321
$ws4.Range("F9").Value = 867
$ws4.Range("F10").Value = 53
$ws4.Range("F11").Value = 517
$ws4.Range("F14").Value = 1144
$ws4.Range("F17").Value = 241
$ws4.Range("F18").Value = 35
$ws4.Range("F19").Value = 414
$ws4.Range("F20").Value = 6646
$ws4.Range("F26").Value = 37
$ws4.Range("F27").Value = 3399
$ws4.Range("F29").Value = 2107
$ws4.Range("F31").Value = 4517
$ws4.Range("F32").Value = 143
$ws4.Range("F38").Value = 1694
$ws4.Range("F40").Value = 170
$ws4.Range("F44").Value = 1208
$ws4.Range("F45").Value = 1783
$ws4.Range("F47").Value = 2138
$ws4.Range("F49").Value = 81
